# excel2skos-holocaust_geographies2.xlsx
#
# "Small tweaks to URIs":
#  - URIs for narrower concepts are no longer built hierarchically (they
#    could, in theory, have more than one broader term), so every URI now
#    sits directly under the scheme's base URI instead of nesting under its
#    broader concept's path segment.
#  - the base URI moves from the old portal (https://portal.ehri-project.eu)
#    to the data host (http://data.ehri-project.eu), and identifiers switch
#    from hyphens to underscores (ehri-holocaust-geographies ->
#    ehri_holocaust_geographies).
#  - the hyperlinks littered across the "skos:broader" (D) column / header
#    mailto-links are pruned back down to just the four header-label
#    hyperlinks, with qualified (dct:/skos:) display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

# --- 1. Rewrite the ConceptScheme URI + all concept/broader URIs ---------

$ws.Range('B1').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies'

$ws.Range('A8').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/places_of_persecution'

$ws.Range('A9').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/ghettos'
$ws.Range('D9').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/places_of_persecution'

$ws.Range('A10').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/camps'
$ws.Range('D10').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/places_of_persecution'

$ws.Range('A11').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/other_extermination'
$ws.Range('D11').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/places_of_persecution'

$ws.Range('A12').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/sites_violence'
$ws.Range('D12').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/places_of_persecution'

$ws.Range('A13').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/forced_labour'
$ws.Range('D13').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/places_of_persecution'

$ws.Range('A14').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/borders'

$ws.Range('A15').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/borders/state'
$ws.Range('D15').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/borders'

$ws.Range('A16').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/borders/administrative'
$ws.Range('D16').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/borders'

$ws.Range('A17').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A18').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/community'
$ws.Range('D18').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A19').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/residence'
$ws.Range('D19').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A20').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/religious'
$ws.Range('D20').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A21').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/memorial'
$ws.Range('D21').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A22').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/hiding'
$ws.Range('D22').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A23').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/incidents'
$ws.Range('D23').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/local'

$ws.Range('A24').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/forced_mobility'

$ws.Range('A25').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/assembly_points'
$ws.Range('D25').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/forced_mobility'

$ws.Range('A26').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/transport_lines'
$ws.Range('D26').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/forced_mobility'

$ws.Range('A27').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/deportations'
$ws.Range('D27').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/forced_mobility'

$ws.Range('A28').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/death_marches'
$ws.Range('D28').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/forced_mobility'

$ws.Range('A29').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_spaces'

$ws.Range('A30').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_camps'
$ws.Range('D30').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_spaces'

$ws.Range('A31').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/border_crossings'
$ws.Range('D31').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_spaces'

$ws.Range('A32').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/boat_passages'
$ws.Range('D32').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_spaces'

$ws.Range('A33').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/consulates'
$ws.Range('D33').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_spaces'

$ws.Range('A34').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/aid_organisations'
$ws.Range('D34').Value = 'http://data.ehri-project.eu/vocabularies/ehri_holocaust_geographies/refugee_spaces'

# --- 2. Hyperlinks: keep only the 4 header mailto-links, with qualified ---
#        (dct:/skos:) display text; drop all the URI hyperlinks that used
#        to shadow columns A and D on rows 8-34.
#        (Hyperlinks.Delete() on any range wipes every hyperlink on the
#        sheet, so clear them all and re-add just the ones we want.)

$ws.Range('A1').Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range('A2'), 'mailto:title@en', '', '', 'dct:title@en')
$ws.Hyperlinks.Add($ws.Range('A3'), 'mailto:description@en', '', '', 'dct:description@en')
$ws.Hyperlinks.Add($ws.Range('B7'), 'mailto:prefLabel@en', '', '', 'skos:prefLabel@en')
$ws.Hyperlinks.Add($ws.Range('C7'), 'mailto:definition@en', '', '', 'skos:definition@en')

# --- 3. Minor layout tweaks that came along with the edit -----------------

$ws.Rows.Item(18).RowHeight = 14.9
$ws.Columns.Item(1).ColumnWidth = 111.63

$ws.Range('A12').Select() | Out-Null
